$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------
# E2/E3: "ZCTP" / "A003" -> "Z000"
$ws.Range("E2").Value = "Z000"
$ws.Range("E3").Value = "Z000"

# P2/Q2/P3/Q3: updated order numbers
$ws.Range("P2").Value = 4600244281
$ws.Range("Q2").Value = 4503342051
$ws.Range("P3").Value = 4600244282
$ws.Range("Q3").Value = 4503342052

# --- O2 style: underline font + centered date number format --------
$o2 = $ws.Range("O2")
$o2.Font.Underline = 2
$o2.NumberFormatLocal = "mm-dd-yy"
$o2.HorizontalAlignment = -4108

# --- Column E width (best-fit widened) ------------------------------
$ws.Columns("E:E").ColumnWidth = 10.71

# --- Sheet view: scroll back to A1, select P2:Q3 --------------------
$win = $excel.Windows.Item(1)
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("P2:Q3").Select() | Out-Null

# --- Page setup -------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
